# Combine testing and the live branch for adding firefighter, station user, and vehicle user.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Set plain (non-hyperlinked) values first, in the column-major order the
#     original author entered them, so the shared-string table comes out in
#     the same sequence as the source edit. ---

# Row 13
$ws.Range("A13").Value = "St1AttributeName"
$ws.Range("A14").Value = "St1ResourceName"

$ws.Range("B13").Value = "ADST06A2"
$ws.Range("B14").Value = "ADST06R2"

$ws.Range("C13").Value = "D1ST4A1"
$ws.Range("C14").Value = "D1ST04R1"

$ws.Range("C15").Value = "bhupesh+TestingStationUser1@atinatechnology.in"
$ws.Range("C16").Value = "bhupesh+TestingFirefighter1@atinatechnology.in"
$ws.Range("C17").Value = "bhupesh+TestingFirefighter2@atinatechnology.in"

$ws.Range("A15").Value = "stUEmail"
$ws.Range("A16").Value = "ff1UEmail"
$ws.Range("A17").Value = "ff2UEmail"

$ws.Range("B15").Value = "bhupesh+LiveStationUser1@atinatechnology.in"
$ws.Range("B16").Value = "bhupesh+LiveFirefighter1@atinatechnology.in"
$ws.Range("B17").Value = "bhupesh+LiveFirefighter2@atinatechnology.in"

$ws.Range("C11").Value = "Bhupesh+d1St04V1newDaily@atinatechnology.in"
$ws.Range("C12").Value = "Bhupesh+d1St04V2newDaily@atinatechnology.in"

# --- Now add the hyperlinks (values already in place, so Add() just wires
#     the relationship without re-touching the shared string table), in the
#     order the author clicked "Insert Hyperlink" on each cell. ---

$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:bhupesh+TestingStationUser1@atinatechnology.in") | Out-Null
$ws.Range("C15").Value = "bhupesh+TestingStationUser1@atinatechnology.in"

$ws.Hyperlinks.Add($ws.Range("C17"), "mailto:bhupesh+TestingFirefighter2@atinatechnology.in") | Out-Null
$ws.Range("C17").Value = "bhupesh+TestingFirefighter2@atinatechnology.in"

$ws.Hyperlinks.Add($ws.Range("B15"), "mailto:bhupesh+LiveStationUser1@atinatechnology.in") | Out-Null
$ws.Range("B15").Value = "bhupesh+LiveStationUser1@atinatechnology.in"

$ws.Hyperlinks.Add($ws.Range("B17"), "mailto:bhupesh+LiveFirefighter2@atinatechnology.in") | Out-Null
$ws.Range("B17").Value = "bhupesh+LiveFirefighter2@atinatechnology.in"

$ws.Hyperlinks.Add($ws.Range("B16"), "mailto:bhupesh+LiveFirefighter1@atinatechnology.in") | Out-Null
$ws.Range("B16").Value = "bhupesh+LiveFirefighter1@atinatechnology.in"

$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:Bhupesh+d1St04V1newDaily@atinatechnology.in") | Out-Null
$ws.Range("C11").Value = "Bhupesh+d1St04V1newDaily@atinatechnology.in"

$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:Bhupesh+d1St04V2newDaily@atinatechnology.in") | Out-Null
$ws.Range("C12").Value = "Bhupesh+d1St04V2newDaily@atinatechnology.in"

# --- Column B widened to fit the new, longer email addresses ---
$ws.Columns.Item(2).ColumnWidth = 48.5

# --- Selection left on C12 by the author when they finished editing ---
$ws.Range("C12").Select() | Out-Null
